# Update cryptos list price/volume data ($excel.ActiveWorkbook is already open)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.548.48'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.563.83'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.988'
$ws.Range("E4").Value = '  -1.72%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.63'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  -1.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.51'
$ws.Range("E8").Value = '  +1.96%  '
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("E10").Value = '  -0.18%  '
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.786.52'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.547.64'
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.500.67'
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.49'
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '224.57'
$ws.Range("E18").Value = '  +4.03%  '
$ws.Range("E19").Value = '  +1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0706'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("E21").Value = '  -1.68%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.03'
$ws.Range("E25").Value = '  -2.05%  '
$ws.Range("E26").Value = '  +2.21%  '
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("E28").Value = '  +0.78%  '
$ws.Range("E29").Value = '  -1.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.14'
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("E32").Value = '  +0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.453.28'
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.16'
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("E35").Value = '  +3.13%  '
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("E39").Value = '  +1.54%  '
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("E41").Value = '  -1.52%  '
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("E44").Value = '  +7.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.968'
$ws.Range("E45").Value = '  -3.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.89'
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.700.05'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.47'
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0523'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₇0976'
$ws.Range("E50").Value = '  -6.00%  '
$ws.Range("E51").Value = '  -1.13%  '